$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snapshot")

# New "dias" (days) values for column B, rows 2-30
$newValues = @(12156, 12154, 12152, 12153, 12151, 12079, 12119, 12116, 12104, 12086, 12092, 12076, 12073, 12074, 11977, 11984, 11956, 10137, 11906, 11740, 11462, 11449, 11620, 11339, 1421, 7735, 1097, 1271, 8031)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Remove the now-unused trailing rows (old rows 31-36)
$ws.Range("A31:D36").Delete(-4162) | Out-Null
